$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Control-hazard fix: B7 ("xx") becomes "11" (Jump row no longer don't-care),
# matching the red highlight style already used for the "Jump" header (F1).
$ws.Range("B7").Value = "11"
$ws.Range("B7").Font.Color = $ws.Range("F1").Font.Color

# Move the active selection from E8 to E10
$ws.Range("E10").Select()
